$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both contain the same data table; update the
# "想去人数" (want-to-go count) values for rows 2 and 3 in each.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 323
    $ws.Range("F3").Value = 1333
}
